$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row: _old -> _FV2404, _new -> _FV2410
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value() -replace "_old$", "_FV2404")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value() -replace "_new$", "_FV2410")
}

# Convert the range into an Excel Table (ListObject)
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U85"), [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Freeze the header row (pane split after row 1)
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
